# MasterDegree_StudyPlan.xlsx update:
#  - "Cyber physical systems and IoT security" (row 18) and
#    "Wireless networks for mobile applications" (row 17) swap places:
#    the former moves up to row 17 and is marked as PASSED (and loses the
#    "in progress" yellow highlight); the latter moves down to row 18 and
#    keeps the yellow "in progress" highlight.
#  - Row 15 ("Law and data") loses its (already invisible) highlight flag.
#  - The credit note on "Internet of things and smart cities" (I19) is
#    updated from "5/6 credits awarded" to
#    "4/6 credits awarded (waiting for last 2 CFU)" and wraps.
#  - The view scrolls down a bit and the active cell moves to H21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Swap the content+formatting of row 17 and row 18 (columns A:H),
#    using a scratch row far below the data as a temporary buffer so
#    nothing is lost mid-swap.
# ---------------------------------------------------------------------
$bufferRow = 500
$ws.Range("A17:H17").Copy($ws.Range("A" + $bufferRow + ":H" + $bufferRow))
$ws.Range("A18:H18").Copy($ws.Range("A17:H17"))
$ws.Range("A" + $bufferRow + ":H" + $bufferRow).Copy($ws.Range("A18:H18"))
$ws.Range("A" + $bufferRow + ":H" + $bufferRow).Clear()

# ---------------------------------------------------------------------
# 2. Row 17 now holds "Cyber physical systems and IoT security": the
#    course is finished, so remove the "in progress" yellow highlight
#    and mark the PASSED column with an "X".
# ---------------------------------------------------------------------
$ws.Range("A17:H17").Interior.ColorIndex = -4142
$ws.Range("H17").Value = "X"

# ---------------------------------------------------------------------
# 3. Row 15 ("Law and data") no longer needs its highlight styling
#    either.
# ---------------------------------------------------------------------
$ws.Range("A15:H15").Interior.ColorIndex = -4142

# ---------------------------------------------------------------------
# 4. Update the partial-credit note for "Internet of things and smart
#    cities" in column I of row 19.
# ---------------------------------------------------------------------
$ws.Range("I19").Value = "4/6 credits awarded (waiting for last 2 CFU)"
$ws.Range("I19").WrapText = $true

# ---------------------------------------------------------------------
# 5. Update the saved view state: scroll a bit further down and leave
#    the active cell on H21.
# ---------------------------------------------------------------------
$ws.Range("A17").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H21").Select()
